$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 100
$ws.Range("E5").Value = 100
$ws.Range("E6").Value = 100
$ws.Range("E7").Value = 100
$ws.Range("E8").Value = 100
$ws.Range("E9").Value = 100
$ws.Range("E10").Value = 100
$ws.Range("E11").Value = 100
$ws.Range("E15").Value = 100
$ws.Range("E21").Value = 100
$ws.Range("E36").Value = 100
